# "Fruta / hortaliza, semanal" update:
# A new weekly price-report row is inserted at row 590 (pushing the
# existing rows 590-670 down to 591-671), growing the used range from
# A1:R670 to A1:R671.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 590; everything below shifts down one row.
$ws.Rows(590).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A590").Value = 5
$ws.Range("B590").Value = "Macroferia Regional de Talca"
$ws.Range("C590").Value = "Maule"
$ws.Range("D590").Value = 45077
$ws.Range("E590").Value = 7
$ws.Range("F590").Value = 100112043
$ws.Range("G590").Value = "Pepino ensalada"
$ws.Range("H590").Value = "Sin especificar"
$ws.Range("I590").Value = "Primera"
$ws.Range("J590").Value = 400
$ws.Range("K590").Value = 10000
$ws.Range("L590").Value = 10000
$ws.Range("M590").Value = 10000
$ws.Range("N590").Value = "`$/caja 60 unidades"
$ws.Range("O590").Value = "Región de Arica y Parinacota"
$ws.Range("P590").Value = 167
$ws.Range("Q590").Value = 60
$ws.Range("R590").Value = "Hortaliza"
